# "added doc, edited maze"
# Adds a new row to the maze-resources table: a YouTube video about the
# recursive-backtracker maze generation algorithm, added by Kaw, and
# widens column B a bit so the longer title fits better. Also updates the
# current selection/scroll position on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Widen column B ------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 42.67

# --- Append the new row of data (row 7) ----------------------------------
$url   = "https://www.youtube.com/watch?v=elMXlO28Q1U&t=265s"
$title = "Maze Generation Algorithm - Recursive Backtracker"
$use   = "comprendre recursive backtracking"
$who   = "Kaw"

$ws.Range("A7").Value = $who
# Set the hyperlink's text/value first so the URL string is registered in
# the shared-strings table before the title text.
$ws.Range("D7").Value = $url
$ws.Range("B7").Value = $title
$ws.Range("C7").Value = $use

# Turn D7 into a clickable hyperlink, matching the style of D2:D6 above it.
$ws.Hyperlinks.Add($ws.Range("D7"), $url)
$ws.Range("D7").Style = $ws.Range("D6").Style

# Match row height/formatting of the other data rows.
$ws.Rows.Item(7).RowHeight = $ws.Rows.Item(6).RowHeight

# --- Update view state -----------------------------------------------
$ws.Range("C6").Select()
$excel.ActiveWindow.ScrollColumn = 2
